# New crime data collected - update weekly CompStat figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header text (Volume/Number and date range) ---
# "Volume 31   Number  41" -> "Volume 31   Number  42"
$ws.Range("A8").Value = "Volume 31   Number  42"

# "Report Covering the Week  10/7/2024  Through  10/13/2024"
# -> "Report Covering the Week  10/14/2024  Through  10/20/2024"
$ws.Range("C9").Value = "Report Covering the Week  10/14/2024  Through  10/20/2024"

# --- Update weekly crime statistics table (rows 14-33) ---

    # Row 14
    $ws.Range("F14").Value = 1
    $ws.Range("G14").Value = 3
    $ws.Range("H14").Value = -66.666666666666
    $ws.Range("M14").Value = -2.272727272727
    $ws.Range("N14").Value = -83.955223880597

    # Row 15
    $ws.Range("C15").Value = 6
    $ws.Range("D15").Value = 4
    $ws.Range("E15").Value = 50
    $ws.Range("F15").Value = 18
    $ws.Range("G15").Value = 10
    $ws.Range("H15").Value = 80
    $ws.Range("I15").Value = 124
    $ws.Range("J15").Value = 108
    $ws.Range("K15").Value = 14.814814814814
    $ws.Range("L15").Value = -18.421052631578
    $ws.Range("M15").Value = -26.190476190476
    $ws.Range("N15").Value = -70.192307692307

    # Row 16
    $ws.Range("C16").Value = 38
    $ws.Range("D16").Value = 47
    $ws.Range("E16").Value = -19.148936170212
    $ws.Range("F16").Value = 152
    $ws.Range("G16").Value = 160
    $ws.Range("H16").Value = -5
    $ws.Range("I16").Value = 1644
    $ws.Range("J16").Value = 1593
    $ws.Range("K16").Value = 3.201506591337
    $ws.Range("L16").Value = -4.139941690962
    $ws.Range("M16").Value = -12.969825304393
    $ws.Range("N16").Value = -78.800773694390

    # Row 17
    $ws.Range("C17").Value = 55
    $ws.Range("D17").Value = 53
    $ws.Range("E17").Value = 3.773584905660
    $ws.Range("F17").Value = 252
    $ws.Range("G17").Value = 238
    $ws.Range("H17").Value = 5.882352941176
    $ws.Range("I17").Value = 2623
    $ws.Range("J17").Value = 2467
    $ws.Range("K17").Value = 6.323469801378
    $ws.Range("L17").Value = 6.669377795851
    $ws.Range("M17").Value = 66.751430387794
    $ws.Range("N17").Value = -45.331388078366

    # Row 18
    $ws.Range("C18").Value = 25
    $ws.Range("D18").Value = 31
    $ws.Range("E18").Value = -19.354838709677
    $ws.Range("F18").Value = 112
    $ws.Range("G18").Value = 109
    $ws.Range("H18").Value = 2.752293577981
    $ws.Range("I18").Value = 1082
    $ws.Range("J18").Value = 1240
    $ws.Range("K18").Value = -12.741935483871
    $ws.Range("L18").Value = -26.344452008168
    $ws.Range("M18").Value = -5.170902716914
    $ws.Range("N18").Value = -88.104661389621

    # Row 19
    $ws.Range("C19").Value = 129
    $ws.Range("D19").Value = 135
    $ws.Range("E19").Value = -4.444444444444
    $ws.Range("F19").Value = 498
    $ws.Range("G19").Value = 506
    $ws.Range("H19").Value = -1.581027667984
    $ws.Range("I19").Value = 5176
    $ws.Range("J19").Value = 5267
    $ws.Range("K19").Value = -1.727738750711
    $ws.Range("L19").Value = -2.725051682014
    $ws.Range("M19").Value = 38.174052322477
    $ws.Range("N19").Value = -41.842696629213

    # Row 20
    $ws.Range("C20").Value = 20
    $ws.Range("D20").Value = 17
    $ws.Range("E20").Value = 17.647058823529
    $ws.Range("G20").Value = 96
    $ws.Range("H20").Value = -13.541666666666
    $ws.Range("I20").Value = 779
    $ws.Range("J20").Value = 1089
    $ws.Range("K20").Value = -28.466483011937
    $ws.Range("L20").Value = -24.442289039767
    $ws.Range("M20").Value = 62.970711297071
    $ws.Range("N20").Value = -89.812998561527

    # Row 21
    $ws.Range("C21").Value = 273
    $ws.Range("D21").Value = 287
    $ws.Range("E21").Value = -4.878048780487
    $ws.Range("F21").Value = 1116
    $ws.Range("G21").Value = 1122
    $ws.Range("H21").Value = -0.534759358288
    $ws.Range("I21").Value = 11471
    $ws.Range("J21").Value = 11811
    $ws.Range("K21").Value = -2.878672424011
    $ws.Range("L21").Value = -5.828749692143
    $ws.Range("M21").Value = 26.905631153888
    $ws.Range("N21").Value = -70.496399176954

    # Row 22
    $ws.Range("C22").Value = 7
    $ws.Range("D22").Value = 16
    $ws.Range("E22").Value = -56.25
    $ws.Range("G22").Value = 42
    $ws.Range("H22").Value = -47.619047619047
    $ws.Range("I22").Value = 207
    $ws.Range("J22").Value = 253
    $ws.Range("K22").Value = -18.181818181818
    $ws.Range("L22").Value = -17.2
    $ws.Range("M22").Value = 10.106382978723

    # Row 23
    $ws.Range("C23").Value = 21
    $ws.Range("D23").Value = 30
    $ws.Range("E23").Value = -30
    $ws.Range("F23").Value = 102
    $ws.Range("G23").Value = 105
    $ws.Range("H23").Value = -2.857142857142
    $ws.Range("I23").Value = 1058
    $ws.Range("J23").Value = 1050
    $ws.Range("K23").Value = 0.761904761904
    $ws.Range("L23").Value = 1.244019138755
    $ws.Range("M23").Value = 54.452554744525

    # Row 24
    $ws.Range("C24").Value = 329
    $ws.Range("D24").Value = 241
    $ws.Range("E24").Value = 36.514522821576
    $ws.Range("F24").Value = 1192
    $ws.Range("G24").Value = 987
    $ws.Range("H24").Value = 20.770010131712
    $ws.Range("I24").Value = 10762
    $ws.Range("J24").Value = 11253
    $ws.Range("K24").Value = -4.363280902870
    $ws.Range("L24").Value = -16.275089466314
    $ws.Range("M24").Value = 39.875227449961

    # Row 25
    $ws.Range("C25").Value = 186
    $ws.Range("D25").Value = 121
    $ws.Range("E25").Value = 53.719008264462
    $ws.Range("F25").Value = 690
    $ws.Range("G25").Value = 517
    $ws.Range("H25").Value = 33.462282398452
    $ws.Range("I25").Value = 5798
    $ws.Range("J25").Value = 6271
    $ws.Range("K25").Value = -7.542656673576
    $ws.Range("L25").Value = -28.736479842674

    # Row 26
    $ws.Range("C26").Value = 102
    $ws.Range("D26").Value = 97
    $ws.Range("E26").Value = 5.154639175257
    $ws.Range("G26").Value = 380
    $ws.Range("H26").Value = 4.473684210526
    $ws.Range("I26").Value = 4185
    $ws.Range("J26").Value = 3697
    $ws.Range("K26").Value = 13.199891804165
    $ws.Range("L26").Value = 15.575807787903
    $ws.Range("M26").Value = -4.233409610983

    # Row 27
    $ws.Range("C27").Value = 7
    $ws.Range("D27").Value = 7
    $ws.Range("F27").Value = 25
    $ws.Range("G27").Value = 19
    $ws.Range("H27").Value = 31.578947368421
    $ws.Range("I27").Value = 195
    $ws.Range("J27").Value = 203
    $ws.Range("K27").Value = -3.940886699507
    $ws.Range("L27").Value = -19.087136929460

    # Row 28
    $ws.Range("C28").Value = 10
    $ws.Range("D28").Value = 7
    $ws.Range("E28").Value = 42.857142857142
    $ws.Range("F28").Value = 50
    $ws.Range("G28").Value = 44
    $ws.Range("H28").Value = 13.636363636363
    $ws.Range("I28").Value = 488
    $ws.Range("J28").Value = 472
    $ws.Range("K28").Value = 3.389830508474
    $ws.Range("L28").Value = -10.622710622710

    # Row 29
    $ws.Range("D29").Value = 2
    $ws.Range("E29").Value = 0
    $ws.Range("F29").Value = 10
    $ws.Range("G29").Value = 10
    $ws.Range("H29").Value = 0
    $ws.Range("I29").Value = 105
    $ws.Range("J29").Value = 126
    $ws.Range("K29").Value = -16.666666666666
    $ws.Range("L29").Value = -38.235294117647
    $ws.Range("M29").Value = -34.375
    $ws.Range("N29").Value = -84.066767830045

    # Row 30
    $ws.Range("F30").Value = 9
    $ws.Range("G30").Value = 8
    $ws.Range("H30").Value = 12.5
    $ws.Range("I30").Value = 85
    $ws.Range("J30").Value = 112
    $ws.Range("K30").Value = -24.107142857142
    $ws.Range("L30").Value = -39.285714285714
    $ws.Range("M30").Value = -38.848920863309
    $ws.Range("N30").Value = -85.833333333333

    # Row 31
    $ws.Range("D31").Value = 7
    $ws.Range("F31").Value = 6
    $ws.Range("G31").Value = 18
    $ws.Range("H31").Value = -66.666666666666
    $ws.Range("I31").Value = 91
    $ws.Range("J31").Value = 78
    $ws.Range("K31").Value = 16.666666666666
    $ws.Range("L31").Value = 13.75

    # Row 33
    $ws.Range("E33").Value = -100
    $ws.Range("J33").Value = 19
    $ws.Range("K33").Value = 10.526315789473

# --- Row 33 "Traffic Fatalities": Week-to-date 2024 count becomes 0, stored
#     as the text "0" (shared string), matching the convention already used
#     elsewhere in the sheet (e.g. C14/D14) for zero counts in this column.
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C33").PasteSpecial(-4122)  # xlPasteFormats - restore General style (s=13)
